$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Report Product Low Stock"
$ws.Range("F1").Value = "In-Stock"
$ws.Range("F2").Select()
